$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.139.89"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "2.929.71"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.05"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.24"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.00"
$ws.Range("E9").Value = "  +4.09%  "
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000225"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.75"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "3.414.74"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").Value = "61.135.63"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").Value = "2.930.99"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "434.71"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.46"
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.679"
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.50"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.09"
$ws.Range("E24").Value = "  +2.94%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.85"
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.96"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.110"
$ws.Range("E31").Value = "  +2.77%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.77"
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "0.0₃0870"
$ws.Range("E34").Value = "  +2.49%  "
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.98"
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "42.05"
$ws.Range("E41").Value = "  +4.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "372.95"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0346"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").Value = "2.709.75"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.00"
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.83"
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.105"
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.00"
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("E51").Value = "  -0.24%  "
